$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B" = 0.6383931775788736
    "C" = -23.10748969621655
    "D" = -0.1173543786648419
    "E" = 0.5034133321369335
    "F" = 0.323629019790324
    "G" = 0.2146654303905574
    "H" = 14.31124727300474
    "I" = 0.2655563776096687
    "J" = 0.2554016671446476
    "K" = 0.2604790223771581
    "L" = 0.2718993657310753
    "M" = 0.4633200086231518
    "N" = -0.08482046726337922
    "O" = 0.4830445005948765
    "P" = 35.07734920122535
    "Q" = 54.57936239911656
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
